$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Replace the old "Participant ID" query in B2 with the new, corrected
# query (the failed test case fix described in the commit message).
$newParticipantQuery = "MATCH (p:participant)-->(s:study)`nOPTIONAL MATCH (samp:sample)-->(p)`nOPTIONAL MATCH (p)<--(diag:diagnosis)`nOPTIONAL MATCH (samp)<--(f:file)`nOPTIONAL MATCH (f)<--(g:genomic_info)`nWITH s, p, samp, f, g, diag`nWHERE g.instrument_model in ['DNBSEQ-G400']`nWITH p`nOPTIONAL MATCH (p)-->(s:study)`nOPTIONAL MATCH (samp:sample)-->(p)`nWITH s, p, apoc.coll.sort(collect(distinct samp.sample_id)) as samp`nRETURN `ncoalesce(p.participant_id,'') as ``Participant ID``,`ncoalesce(s.study_name, '') as ``Study Name``,`ncoalesce(s.phs_accession,'') as ``Accession``,`ncoalesce(p.gender,'') as ``Gender``,`ncoalesce(apoc.text.join(samp, ','), '') as ``Samples```nORDER BY p.participant_id limit 100"

$ws.Range("B2").Value = $newParticipantQuery

# Update the active selection to match the saved view state (C2).
$ws.Range("C2").Select()
